{"js": "// Add an \"OBJECTIVE:\" section (styled like the other resume section\n// headings) right before the \"EDUCATION:\" heading, and move the\n// \"_GoBack\" bookmark (an artifact Word leaves at the last edited\n// location) from the end of the document to the point inside the new\n// paragraph where the author's typing paused.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\n// Bookmark names must be unique, so drop the existing \"_GoBack\" before\n// re-creating it further down at its new location.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Locate the \"EDUCATION:\" heading paragraph - the new OBJECTIVE\n// paragraph goes directly above it.\nlet educationParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"EDUCATION:\") {\n    educationParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!educationParagraph) {\n  throw new Error('Could not find the \"EDUCATION:\" heading paragraph.');\n}\n\n// Create the new heading-styled paragraph with its first run.\nconst objectiveParagraph = educationParagraph.insertParagraph(\n  \"OBJECTIVE: \",\n  Word.InsertLocation.before\n);\nobjectiveParagraph.styleBuiltIn = Word.BuiltInStyleName.heading1;\n\n// \"Eager to drive\" is typed next, at a smaller (12pt) size than the\n// heading default.\nconst eagerRun = objectiveParagraph.insertText(\n  \"Eager to drive\",\n  Word.InsertLocation.end\n);\neagerRun.font.size = 12;\nawait context.sync();\n\n// This is where the cursor was left last - Word records that as the\n// hidden \"_GoBack\" bookmark.\nconst goBackRange = eagerRun.getRange(Word.RangeLocation.end);\ngoBackRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n\nconst solutionsAtRun = objectiveParagraph.insertText(\n  \" solutions at \",\n  Word.InsertLocation.end\n);\nsolutionsAtRun.font.size = 12;\n\nconst acuteEngineeringRun = objectiveParagraph.insertText(\n  \"Acute Engineering\",\n  Word.InsertLocation.end\n);\nacuteEngineeringRun.font.size = 12;\n\nconst fullTimeRun = objectiveParagraph.insertText(\n  \" on a full-time basis\",\n  Word.InsertLocation.end\n);\nfullTimeRun.font.size = 12;\n\nawait context.sync();\n", "ps1": "# Add an \"OBJECTIVE:\" section (styled like the other resume section\n# headings) right before the \"EDUCATION:\" heading, and move the\n# \"_GoBack\" bookmark (an artifact Word leaves at the last edited\n# location) from the end of the document to the point inside the new\n# paragraph where the author's typing paused.\n\n$d = $word.ActiveDocument\n\n# Bookmark names must be unique, so drop the existing \"_GoBack\" before\n# re-creating it further down at its new location.\ntry {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n} catch {\n}\n\nfunction Find-ParagraphWithText($doc, $text) {\n    foreach ($p in $doc.Paragraphs) {\n        $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n        if ($t -eq $text) {\n            return $p\n        }\n    }\n    return $null\n}\n\nfunction Find-ParagraphBeforeText($doc, $text) {\n    $previous = $null\n    foreach ($p in $doc.Paragraphs) {\n        $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n        if ($t -eq $text) {\n            return $previous\n        }\n        $previous = $p\n    }\n    return $null\n}\n\n# Insert a new empty paragraph right before \"EDUCATION:\".\n$educationParagraph = Find-ParagraphWithText $d \"EDUCATION:\"\nif ($educationParagraph -eq $null) {\n    throw 'Could not find the \"EDUCATION:\" heading paragraph.'\n}\n$insertionPoint = $educationParagraph.Range\n$insertionPoint.Collapse(1)  # wdCollapseStart\n$insertionPoint.InsertParagraphBefore()\n\n# Style the new (now empty) paragraph like the other section headings\n# and give it its first run of text.\n$objectivePara = Find-ParagraphBeforeText $d \"EDUCATION:\"\n$objectivePara.Range.Style = \"Heading 1\"\n$objectivePara.Range.Text = \"OBJECTIVE: \"\n\nfunction Get-ParagraphEndRange($doc, $para) {\n    # A collapsed range positioned right before the paragraph mark, so\n    # that InsertAfter() appends inside the paragraph instead of\n    # merging with whatever comes next.\n    $r = $doc.Range($para.Range.Start, $para.Range.End - 1)\n    $r.Collapse(0)  # wdCollapseEnd\n    return $r\n}\n\n# \"Eager to drive\" is typed next, at a smaller (12pt) size than the\n# heading default.\n$eagerRange = Get-ParagraphEndRange $d $objectivePara\n$eagerStart = $eagerRange.Start\n$eagerText = \"Eager to drive\"\n$eagerRange.InsertAfter($eagerText)\n$d.Range($eagerStart, $eagerStart + $eagerText.Length).Font.Size = 12\n\n$solutionsAtRange = Get-ParagraphEndRange $d $objectivePara\n$solutionsAtStart = $solutionsAtRange.Start\n$solutionsAtText = \" solutions at \"\n$solutionsAtRange.InsertAfter($solutionsAtText)\n$d.Range($solutionsAtStart, $solutionsAtStart + $solutionsAtText.Length).Font.Size = 12\n\n$acuteRange = Get-ParagraphEndRange $d $objectivePara\n$acuteStart = $acuteRange.Start\n$acuteText = \"Acute Engineering\"\n$acuteRange.InsertAfter($acuteText)\n$d.Range($acuteStart, $acuteStart + $acuteText.Length).Font.Size = 12\n\n$fullTimeRange = Get-ParagraphEndRange $d $objectivePara\n$fullTimeStart = $fullTimeRange.Start\n$fullTimeText = \" on a full-time basis\"\n$fullTimeRange.InsertAfter($fullTimeText)\n$d.Range($fullTimeStart, $fullTimeStart + $fullTimeText.Length).Font.Size = 12\n\n# This is where the cursor was left last (right after \"Eager to\n# drive\") - Word records that as the hidden \"_GoBack\" bookmark. Added\n# last, once the whole paragraph is in place, so the position is mid-\n# paragraph rather than sitting on a paragraph-mark boundary.\n$goBackPos = $eagerStart + $eagerText.Length\n$d.Range($goBackPos, $goBackPos).Bookmarks.Add(\"_GoBack\")\n"}
